$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Apoe"
$ws.Cells.Item(2, 3).Value = "Sorl1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"29.32133366666666"
$ws.Cells.Item(2, 8).Value = [double]"87.964001"
$ws.Cells.Item(2, 9).Value = [double]"0.006401919837078288"
$ws.Cells.Item(2, 10).Value = [double]"0.006401919837078288"
$ws.Cells.Item(2, 11).Value = [double]"2"
$ws.Cells.Item(2, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2, 13).Value = [double]"0.05271666666666667"
$ws.Cells.Item(2, 14).Value = [double]"0.15815"
$ws.Cells.Item(2, 15).Value = [double]"0.003013014833311122"
$ws.Cells.Item(2, 16).Value = [double]"0.003013014833311122"
$ws.Cells.Item(2, 17).Value = [double]"1.545722973127778"
$ws.Cells.Item(2, 18).Value = [double]"13.91150675815"
$ws.Cells.Item(2, 19).Value = [double]"1.92890794307856E-05"
$ws.Cells.Item(2, 20).Value = [double]"1.928907943078561E-05"

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Apoe"
$ws.Cells.Item(3, 3).Value = "Sorl1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"29.32133366666666"
$ws.Cells.Item(3, 8).Value = [double]"87.964001"
$ws.Cells.Item(3, 9).Value = [double]"0.006401919837078288"
$ws.Cells.Item(3, 10).Value = [double]"0.006401919837078288"
$ws.Cells.Item(3, 11).Value = [double]"3"
$ws.Cells.Item(3, 12).Value = [double]"1"
$ws.Cells.Item(3, 13).Value = [double]"1.046984666666667"
$ws.Cells.Item(3, 14).Value = [double]"3.140954"
$ws.Cells.Item(3, 15).Value = [double]"0.05984028449413786"
$ws.Cells.Item(3, 16).Value = [double]"0.05984028449413786"
$ws.Cells.Item(3, 17).Value = [double]"30.69898675521711"
$ws.Cells.Item(3, 18).Value = [double]"276.290880796954"
$ws.Cells.Item(3, 19).Value = [double]"0.0003830927043594294"
$ws.Cells.Item(3, 20).Value = [double]"0.0003830927043594295"

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Apoe"
$ws.Cells.Item(4, 3).Value = "Sorl1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"29.32133366666666"
$ws.Cells.Item(4, 8).Value = [double]"87.964001"
$ws.Cells.Item(4, 9).Value = [double]"0.006401919837078288"
$ws.Cells.Item(4, 10).Value = [double]"0.006401919837078288"
$ws.Cells.Item(4, 11).Value = [double]"3"
$ws.Cells.Item(4, 12).Value = [double]"1"
$ws.Cells.Item(4, 13).Value = [double]"15.25299333333333"
$ws.Cells.Item(4, 14).Value = [double]"45.75898"
$ws.Cells.Item(4, 15).Value = [double]"0.8717830255908123"
$ws.Cells.Item(4, 16).Value = [double]"0.8717830255908123"
$ws.Cells.Item(4, 17).Value = [double]"447.2381069421089"
$ws.Cells.Item(4, 18).Value = [double]"4025.14296247898"
$ws.Cells.Item(4, 19).Value = [double]"0.00558108504515795"
$ws.Cells.Item(4, 20).Value = [double]"0.005581085045157951"

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Apoe"
$ws.Cells.Item(5, 3).Value = "Sorl1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"29.32133366666666"
$ws.Cells.Item(5, 8).Value = [double]"87.964001"
$ws.Cells.Item(5, 9).Value = [double]"0.006401919837078288"
$ws.Cells.Item(5, 10).Value = [double]"0.006401919837078288"
$ws.Cells.Item(5, 11).Value = [double]"3"
$ws.Cells.Item(5, 12).Value = [double]"1"
$ws.Cells.Item(5, 13).Value = [double]"1.143623666666667"
$ws.Cells.Item(5, 14).Value = [double]"3.430871"
$ws.Cells.Item(5, 15).Value = [double]"0.06536367508173863"
$ws.Cells.Item(5, 16).Value = [double]"0.06536367508173863"
$ws.Cells.Item(5, 17).Value = [double]"33.53257111943011"
$ws.Cells.Item(5, 18).Value = [double]"301.793140074871"
$ws.Cells.Item(5, 19).Value = [double]"0.0004184530081301223"
$ws.Cells.Item(5, 20).Value = [double]"0.0004184530081301223"

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Apoe"
$ws.Cells.Item(6, 3).Value = "Sorl1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = [double]"3"
$ws.Cells.Item(6, 6).Value = [double]"1"
$ws.Cells.Item(6, 7).Value = [double]"45.524413"
$ws.Cells.Item(6, 8).Value = [double]"136.573239"
$ws.Cells.Item(6, 9).Value = [double]"0.009939644832300594"
$ws.Cells.Item(6, 10).Value = [double]"0.009939644832300592"
$ws.Cells.Item(6, 11).Value = [double]"2"
$ws.Cells.Item(6, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6, 13).Value = [double]"0.05271666666666667"
$ws.Cells.Item(6, 14).Value = [double]"0.15815"
$ws.Cells.Item(6, 15).Value = [double]"0.003013014833311122"
$ws.Cells.Item(6, 16).Value = [double]"0.003013014833311122"
$ws.Cells.Item(6, 17).Value = [double]"2.399895305316667"
$ws.Cells.Item(6, 18).Value = [double]"21.59905774785"
$ws.Cells.Item(6, 19).Value = [double]"2.994829731756593E-05"
$ws.Cells.Item(6, 20).Value = [double]"2.994829731756593E-05"

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Apoe"
$ws.Cells.Item(7, 3).Value = "Sorl1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = [double]"3"
$ws.Cells.Item(7, 6).Value = [double]"1"
$ws.Cells.Item(7, 7).Value = [double]"45.524413"
$ws.Cells.Item(7, 8).Value = [double]"136.573239"
$ws.Cells.Item(7, 9).Value = [double]"0.009939644832300594"
$ws.Cells.Item(7, 10).Value = [double]"0.009939644832300592"
$ws.Cells.Item(7, 11).Value = [double]"3"
$ws.Cells.Item(7, 12).Value = [double]"1"
$ws.Cells.Item(7, 13).Value = [double]"1.046984666666667"
$ws.Cells.Item(7, 14).Value = [double]"3.140954"
$ws.Cells.Item(7, 15).Value = [double]"0.05984028449413786"
$ws.Cells.Item(7, 16).Value = [double]"0.05984028449413786"
$ws.Cells.Item(7, 17).Value = [double]"47.66336237000067"
$ws.Cells.Item(7, 18).Value = [double]"428.970261330006"
$ws.Cells.Item(7, 19).Value = [double]"0.0005947911745355547"
$ws.Cells.Item(7, 20).Value = [double]"0.0005947911745355546"

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Apoe"
$ws.Cells.Item(8, 3).Value = "Sorl1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = [double]"3"
$ws.Cells.Item(8, 6).Value = [double]"1"
$ws.Cells.Item(8, 7).Value = [double]"45.524413"
$ws.Cells.Item(8, 8).Value = [double]"136.573239"
$ws.Cells.Item(8, 9).Value = [double]"0.009939644832300594"
$ws.Cells.Item(8, 10).Value = [double]"0.009939644832300592"
$ws.Cells.Item(8, 11).Value = [double]"3"
$ws.Cells.Item(8, 12).Value = [double]"1"
$ws.Cells.Item(8, 13).Value = [double]"15.25299333333333"
$ws.Cells.Item(8, 14).Value = [double]"45.75898"
$ws.Cells.Item(8, 15).Value = [double]"0.8717830255908123"
$ws.Cells.Item(8, 16).Value = [double]"0.8717830255908123"
$ws.Cells.Item(8, 17).Value = [double]"694.3835679929134"
$ws.Cells.Item(8, 18).Value = [double]"6249.45211193622"
$ws.Cells.Item(8, 19).Value = [double]"0.008665213645201094"
$ws.Cells.Item(8, 20).Value = [double]"0.008665213645201092"

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Apoe"
$ws.Cells.Item(9, 3).Value = "Sorl1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = [double]"3"
$ws.Cells.Item(9, 6).Value = [double]"1"
$ws.Cells.Item(9, 7).Value = [double]"45.524413"
$ws.Cells.Item(9, 8).Value = [double]"136.573239"
$ws.Cells.Item(9, 9).Value = [double]"0.009939644832300594"
$ws.Cells.Item(9, 10).Value = [double]"0.009939644832300592"
$ws.Cells.Item(9, 11).Value = [double]"3"
$ws.Cells.Item(9, 12).Value = [double]"1"
$ws.Cells.Item(9, 13).Value = [double]"1.143623666666667"
$ws.Cells.Item(9, 14).Value = [double]"3.430871"
$ws.Cells.Item(9, 15).Value = [double]"0.06536367508173863"
$ws.Cells.Item(9, 16).Value = [double]"0.06536367508173863"
$ws.Cells.Item(9, 17).Value = [double]"52.06279611790767"
$ws.Cells.Item(9, 18).Value = [double]"468.565165061169"
$ws.Cells.Item(9, 19).Value = [double]"0.0006496917152463784"
$ws.Cells.Item(9, 20).Value = [double]"0.0006496917152463783"

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Apoe"
$ws.Cells.Item(10, 3).Value = "Sorl1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = [double]"3"
$ws.Cells.Item(10, 6).Value = [double]"1"
$ws.Cells.Item(10, 7).Value = [double]"4438.215250666667"
$ws.Cells.Item(10, 8).Value = [double]"13314.645752"
$ws.Cells.Item(10, 9).Value = [double]"0.9690247577915309"
$ws.Cells.Item(10, 10).Value = [double]"0.9690247577915307"
$ws.Cells.Item(10, 11).Value = [double]"2"
$ws.Cells.Item(10, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10, 13).Value = [double]"0.05271666666666667"
$ws.Cells.Item(10, 14).Value = [double]"0.15815"
$ws.Cells.Item(10, 15).Value = [double]"0.003013014833311122"
$ws.Cells.Item(10, 16).Value = [double]"0.003013014833311122"
$ws.Cells.Item(10, 17).Value = [double]"233.9679139643112"
$ws.Cells.Item(10, 18).Value = [double]"2105.7112256788"
$ws.Cells.Item(10, 19).Value = [double]"0.0029196859690716"
$ws.Cells.Item(10, 20).Value = [double]"0.0029196859690716"

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Apoe"
$ws.Cells.Item(11, 3).Value = "Sorl1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = [double]"3"
$ws.Cells.Item(11, 6).Value = [double]"1"
$ws.Cells.Item(11, 7).Value = [double]"4438.215250666667"
$ws.Cells.Item(11, 8).Value = [double]"13314.645752"
$ws.Cells.Item(11, 9).Value = [double]"0.9690247577915309"
$ws.Cells.Item(11, 10).Value = [double]"0.9690247577915307"
$ws.Cells.Item(11, 11).Value = [double]"3"
$ws.Cells.Item(11, 12).Value = [double]"1"
$ws.Cells.Item(11, 13).Value = [double]"1.046984666666667"
$ws.Cells.Item(11, 14).Value = [double]"3.140954"
$ws.Cells.Item(11, 15).Value = [double]"0.05984028449413786"
$ws.Cells.Item(11, 16).Value = [double]"0.05984028449413786"
$ws.Cells.Item(11, 17).Value = [double]"4646.743314814157"
$ws.Cells.Item(11, 18).Value = [double]"41820.68983332741"
$ws.Cells.Item(11, 19).Value = [double]"0.05798671718810824"
$ws.Cells.Item(11, 20).Value = [double]"0.05798671718810823"

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Apoe"
$ws.Cells.Item(12, 3).Value = "Sorl1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = [double]"3"
$ws.Cells.Item(12, 6).Value = [double]"1"
$ws.Cells.Item(12, 7).Value = [double]"4438.215250666667"
$ws.Cells.Item(12, 8).Value = [double]"13314.645752"
$ws.Cells.Item(12, 9).Value = [double]"0.9690247577915309"
$ws.Cells.Item(12, 10).Value = [double]"0.9690247577915307"
$ws.Cells.Item(12, 11).Value = [double]"3"
$ws.Cells.Item(12, 12).Value = [double]"1"
$ws.Cells.Item(12, 13).Value = [double]"15.25299333333333"
$ws.Cells.Item(12, 14).Value = [double]"45.75898"
$ws.Cells.Item(12, 15).Value = [double]"0.8717830255908123"
$ws.Cells.Item(12, 16).Value = [double]"0.8717830255908123"
$ws.Cells.Item(12, 17).Value = [double]"67696.067630317"
$ws.Cells.Item(12, 18).Value = [double]"609264.608672853"
$ws.Cells.Item(12, 19).Value = [double]"0.8447793352199049"
$ws.Cells.Item(12, 20).Value = [double]"0.8447793352199048"

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Apoe"
$ws.Cells.Item(13, 3).Value = "Sorl1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = [double]"3"
$ws.Cells.Item(13, 6).Value = [double]"1"
$ws.Cells.Item(13, 7).Value = [double]"4438.215250666667"
$ws.Cells.Item(13, 8).Value = [double]"13314.645752"
$ws.Cells.Item(13, 9).Value = [double]"0.9690247577915309"
$ws.Cells.Item(13, 10).Value = [double]"0.9690247577915307"
$ws.Cells.Item(13, 11).Value = [double]"3"
$ws.Cells.Item(13, 12).Value = [double]"1"
$ws.Cells.Item(13, 13).Value = [double]"1.143623666666667"
$ws.Cells.Item(13, 14).Value = [double]"3.430871"
$ws.Cells.Item(13, 15).Value = [double]"0.06536367508173863"
$ws.Cells.Item(13, 16).Value = [double]"0.06536367508173863"
$ws.Cells.Item(13, 17).Value = [double]"5075.647998423334"
$ws.Cells.Item(13, 18).Value = [double]"45680.83198580999"
$ws.Cells.Item(13, 19).Value = [double]"0.0633390194144461"
$ws.Cells.Item(13, 20).Value = [double]"0.06333901941444608"

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Apoe"
$ws.Cells.Item(14, 3).Value = "Sorl1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = [double]"3"
$ws.Cells.Item(14, 6).Value = [double]"1"
$ws.Cells.Item(14, 7).Value = [double]"67.02347933333333"
$ws.Cells.Item(14, 8).Value = [double]"201.070438"
$ws.Cells.Item(14, 9).Value = [double]"0.01463367753909034"
$ws.Cells.Item(14, 10).Value = [double]"0.01463367753909034"
$ws.Cells.Item(14, 11).Value = [double]"2"
$ws.Cells.Item(14, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14, 13).Value = [double]"0.05271666666666667"
$ws.Cells.Item(14, 14).Value = [double]"0.15815"
$ws.Cells.Item(14, 15).Value = [double]"0.003013014833311122"
$ws.Cells.Item(14, 16).Value = [double]"0.003013014833311122"
$ws.Cells.Item(14, 17).Value = [double]"3.533254418855555"
$ws.Cells.Item(14, 18).Value = [double]"31.7992897697"
$ws.Cells.Item(14, 19).Value = [double]"4.409148749117098E-05"
$ws.Cells.Item(14, 20).Value = [double]"4.409148749117099E-05"

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Apoe"
$ws.Cells.Item(15, 3).Value = "Sorl1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = [double]"3"
$ws.Cells.Item(15, 6).Value = [double]"1"
$ws.Cells.Item(15, 7).Value = [double]"67.02347933333333"
$ws.Cells.Item(15, 8).Value = [double]"201.070438"
$ws.Cells.Item(15, 9).Value = [double]"0.01463367753909034"
$ws.Cells.Item(15, 10).Value = [double]"0.01463367753909034"
$ws.Cells.Item(15, 11).Value = [double]"3"
$ws.Cells.Item(15, 12).Value = [double]"1"
$ws.Cells.Item(15, 13).Value = [double]"1.046984666666667"
$ws.Cells.Item(15, 14).Value = [double]"3.140954"
$ws.Cells.Item(15, 15).Value = [double]"0.05984028449413786"
$ws.Cells.Item(15, 16).Value = [double]"0.05984028449413786"
$ws.Cells.Item(15, 17).Value = [double]"70.17255516865022"
$ws.Cells.Item(15, 18).Value = [double]"631.552996517852"
$ws.Cells.Item(15, 19).Value = [double]"0.000875683427134641"
$ws.Cells.Item(15, 20).Value = [double]"0.0008756834271346411"

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Apoe"
$ws.Cells.Item(16, 3).Value = "Sorl1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = [double]"3"
$ws.Cells.Item(16, 6).Value = [double]"1"
$ws.Cells.Item(16, 7).Value = [double]"67.02347933333333"
$ws.Cells.Item(16, 8).Value = [double]"201.070438"
$ws.Cells.Item(16, 9).Value = [double]"0.01463367753909034"
$ws.Cells.Item(16, 10).Value = [double]"0.01463367753909034"
$ws.Cells.Item(16, 11).Value = [double]"3"
$ws.Cells.Item(16, 12).Value = [double]"1"
$ws.Cells.Item(16, 13).Value = [double]"15.25299333333333"
$ws.Cells.Item(16, 14).Value = [double]"45.75898"
$ws.Cells.Item(16, 15).Value = [double]"0.8717830255908123"
$ws.Cells.Item(16, 16).Value = [double]"0.8717830255908123"
$ws.Cells.Item(16, 17).Value = [double]"1022.308683448138"
$ws.Cells.Item(16, 18).Value = [double]"9200.778151033241"
$ws.Cells.Item(16, 19).Value = [double]"0.01275739168054849"
$ws.Cells.Item(16, 20).Value = [double]"0.01275739168054849"

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Apoe"
$ws.Cells.Item(17, 3).Value = "Sorl1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = [double]"3"
$ws.Cells.Item(17, 6).Value = [double]"1"
$ws.Cells.Item(17, 7).Value = [double]"67.02347933333333"
$ws.Cells.Item(17, 8).Value = [double]"201.070438"
$ws.Cells.Item(17, 9).Value = [double]"0.01463367753909034"
$ws.Cells.Item(17, 10).Value = [double]"0.01463367753909034"
$ws.Cells.Item(17, 11).Value = [double]"3"
$ws.Cells.Item(17, 12).Value = [double]"1"
$ws.Cells.Item(17, 13).Value = [double]"1.143623666666667"
$ws.Cells.Item(17, 14).Value = [double]"3.430871"
$ws.Cells.Item(17, 15).Value = [double]"0.06536367508173863"
$ws.Cells.Item(17, 16).Value = [double]"0.06536367508173863"
$ws.Cells.Item(17, 17).Value = [double]"76.64963718794422"
$ws.Cells.Item(17, 18).Value = [double]"689.846734691498"
$ws.Cells.Item(17, 19).Value = [double]"0.0009565109439160373"
$ws.Cells.Item(17, 20).Value = [double]"0.0009565109439160374"
